# Generate Report for handback
# The row for "e4cbc0be-6f0f-4725-9ce9-cf5df27ac965.md" changes its status
# from "Ready for handoff" to "Handed back: in sync with en-US", and the
# "Latest Handback DateTime" is updated on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G3").Value = "2016-01-11 13:03:28"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G3").Value = "2016-01-11 13:05:35"
